$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, border, centered/top alignment) from an existing
# header cell (AC1) onto the new header cells before setting their values.
# -4122 == xlPasteFormats
$ws.Range("AC1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

for ($r = 2; $r -le 42; $r++) {
    $ws.Cells.Item($r, 30).Value = 79
    $ws.Cells.Item($r, 31).Value = 82
    $ws.Cells.Item($r, 32).Value = 0
}
